# Applies the "Updated config file and argument values" commit to the
# Settings sheet of GameGetter/Data/Config.xlsx.
#
# Summary of the change:
#   1. Three config path values are renamed (file moved out of the
#      Epic/Steam subfolders, into a flat Data/Temp folder with a
#      launcher-prefixed name).
#   2. Two brand-new settings rows (SteamGameList / EpicGameList) are
#      inserted above the WBEpicSheet/WBSteamSheet rows, which pushes
#      every row below down by three rows (two new rows of data + one
#      blank separator row, matching the sheet's existing blank-row
#      convention).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Settings")

# 1) Rename three relative-path values in place.
$ws.Range("B13").Value = "Data/Temp/EpicGamesToInstall.xlsx"
$ws.Range("B20").Value = "Data/Temp/SteamGamesToAdd.xlsx"
$ws.Range("B24").Value = "Data/Temp/SteamGamesToInstall.xlsx"

# 2) Insert three new rows above the current row 29 (WBEpicSheet), which
#    shifts WBEpicSheet/WBSteamSheet/EpicCredential/SteamCredential/
#    Steam_Credential/Steam_Failed_Credential/OrchestratorQueueName (and
#    the blank rows between them) down by three rows.
$ws.Range("A29:A31").EntireRow.Insert()

# 3) Populate the two freshly inserted rows with the new settings; the
#    third inserted row (31) stays blank, mirroring the sheet's existing
#    blank-separator-row style.
$ws.Range("A29").Value = "SteamGameList"
$ws.Range("B29").Value = "Data/Input/SteamGameList.xlsx"
$ws.Range("C29").Value = "Relative path for list of free games on steam"

$ws.Range("A30").Value = "EpicGameList"
$ws.Range("B30").Value = "Data/Input/EpicGameList.xlsx"
$ws.Range("C30").Value = "Relative path for list of free games on Epic"

# Move the active selection to mirror the saved workbook state (B26 on
# the Settings sheet was selected when the author saved the file).
$ws.Range("B26").Select()
